# "Tried to implement Penalty Reward System (unfinished)"
# Update the weekly/monthly PO quantity sheets: a value is adjusted and a
# handful of trailing / now-superseded rows are removed.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item(1)

# Requested quantity for the week of row 20 drops from 120 to 58.
$ws1.Cells.Item(20, 2).Value = 58

# Row 21 (week 45137.99999999999 / qty 40) is removed entirely, shifting
# the remaining weeks up by one.
$ws1.Rows.Item(21).Delete()

# The two trailing weeks (formerly rows 26 and 27, now 25 and 26 after the
# shift above) are also removed, leaving data through row 24.
$ws1.Range("A25:B26").EntireRow.Delete()

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item(2)

# Requested quantity for the month of row 8 drops from 160 to 58.
$ws2.Cells.Item(8, 2).Value = 58

# The trailing month (row 12) is removed, leaving data through row 11.
$ws2.Rows.Item(12).Delete()
